# Auto-generated edit script: applies the Tonberry_Profits.xlsx cell-value
# updates (currentAveragePrice / LevePrice / LeveProfit columns) per the
# scheduled-runner diff, sheet by sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 485.6
$ws.Range("I115").Value = 485.6
$ws.Range("K115").Value = 1456.8
$ws.Range("M115").Value = 110.1999999999998
$ws.Range("H138").Value = 1961.5555
$ws.Range("J138").Value = 2497.375
$ws.Range("L138").Value = 7492.125
$ws.Range("N138").Value = -17772.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3179.3538
$ws.Range("I32").Value = 2133.5576
$ws.Range("K32").Value = 2133.5576
$ws.Range("M32").Value = -1846.5576
$ws.Range("H45").Value = 3216828.2
$ws.Range("I45").Value = 9002299
$ws.Range("K45").Value = 9002299
$ws.Range("M45").Value = -9001922
$ws.Range("H61").Value = 2533.0688
$ws.Range("I61").Value = 1675.4762
$ws.Range("K61").Value = 1675.4762
$ws.Range("M61").Value = -1463.4762
$ws.Range("H74").Value = 1139.3928
$ws.Range("I74").Value = 510.55554
$ws.Range("J74").Value = 2271.3
$ws.Range("K74").Value = 510.55554
$ws.Range("L74").Value = 2271.3
$ws.Range("M74").Value = 363.44446
$ws.Range("N74").Value = -4019.3
$ws.Range("H77").Value = 1139.3928
$ws.Range("I77").Value = 510.55554
$ws.Range("J77").Value = 2271.3
$ws.Range("K77").Value = 2552.7777
$ws.Range("L77").Value = 11356.5
$ws.Range("M77").Value = 1815.2223
$ws.Range("N77").Value = -20092.5
$ws.Range("H97").Value = 1615.75
$ws.Range("I97").Value = 1523.5333
$ws.Range("K97").Value = 1523.5333
$ws.Range("M97").Value = -1027.5333
$ws.Range("H109").Value = 61860.715
$ws.Range("J109").Value = 61860.715
$ws.Range("L109").Value = 61860.715
$ws.Range("N109").Value = -64634.715
$ws.Range("H132").Value = 2037.5088
$ws.Range("I132").Value = 1555.7142
$ws.Range("K132").Value = 4667.142599999999
$ws.Range("M132").Value = -2137.142599999999
$ws.Range("H136").Value = 2533.0688
$ws.Range("I136").Value = 1675.4762
$ws.Range("K136").Value = 5026.4286
$ws.Range("M136").Value = -2476.4286

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 40599.8
$ws.Range("I82").Value = 37666.332
$ws.Range("J82").Value = 45000
$ws.Range("K82").Value = 37666.332
$ws.Range("L82").Value = 45000
$ws.Range("M82").Value = -37283.332
$ws.Range("N82").Value = -45766
$ws.Range("H85").Value = 40599.8
$ws.Range("I85").Value = 37666.332
$ws.Range("J85").Value = 45000
$ws.Range("K85").Value = 37666.332
$ws.Range("L85").Value = 45000
$ws.Range("M85").Value = -36340.332
$ws.Range("N85").Value = -47652
$ws.Range("H94").Value = 1630
$ws.Range("I94").Value = 1477.5
$ws.Range("K94").Value = 1477.5
$ws.Range("M94").Value = -1026.5
$ws.Range("H95").Value = 71444.25
$ws.Range("J95").Value = 71444.25
$ws.Range("L95").Value = 71444.25
$ws.Range("N95").Value = -76936.25
$ws.Range("H97").Value = 5229.3335
$ws.Range("I97").Value = 5229.3335
$ws.Range("K97").Value = 5229.3335
$ws.Range("M97").Value = -4238.3335
$ws.Range("H108").Value = 90682.5
$ws.Range("J108").Value = 90682.5
$ws.Range("L108").Value = 90682.5
$ws.Range("N108").Value = -98362.5
$ws.Range("H109").Value = 54000
$ws.Range("J109").Value = 54000
$ws.Range("L109").Value = 54000
$ws.Range("N109").Value = -56774
$ws.Range("H134").Value = 6315.913
$ws.Range("I134").Value = 6713.3
$ws.Range("K134").Value = 20139.9
$ws.Range("M134").Value = -17604.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 43833.332
$ws.Range("J68").Value = 43833.332
$ws.Range("L68").Value = 43833.332
$ws.Range("N68").Value = -45331.332
$ws.Range("H71").Value = 43833.332
$ws.Range("J71").Value = 43833.332
$ws.Range("L71").Value = 131499.996
$ws.Range("N71").Value = -138987.996
$ws.Range("H132").Value = 2428.8333
$ws.Range("I132").Value = 1535
$ws.Range("J132").Value = 4216.5
$ws.Range("K132").Value = 4605
$ws.Range("L132").Value = 12649.5
$ws.Range("M132").Value = -2075
$ws.Range("N132").Value = -17709.5
$ws.Range("H141").Value = 70439
$ws.Range("J141").Value = 70439
$ws.Range("L141").Value = 70439
$ws.Range("N141").Value = -80799

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 49.6
$ws.Range("I12").Value = 21.571428
$ws.Range("J12").Value = 74.125
$ws.Range("K12").Value = 64.71428400000001
$ws.Range("L12").Value = 222.375
$ws.Range("M12").Value = 108.285716
$ws.Range("N12").Value = -568.375
$ws.Range("H132").Value = 790.8570999999999
$ws.Range("J132").Value = 803.41174
$ws.Range("L132").Value = 7230.70566
$ws.Range("N132").Value = -12290.70566
$ws.Range("H140").Value = 3114.5881
$ws.Range("J140").Value = 6098.5713
$ws.Range("L140").Value = 18295.7139
$ws.Range("N140").Value = -28655.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H49").Value = 24999.5
$ws.Range("J49").Value = 24999.5
$ws.Range("L49").Value = 24999.5
$ws.Range("N49").Value = -25367.5
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H74").Value = 21250
$ws.Range("J74").Value = 21250
$ws.Range("L74").Value = 21250
$ws.Range("N74").Value = -23122
$ws.Range("H77").Value = 21250
$ws.Range("J77").Value = 21250
$ws.Range("L77").Value = 63750
$ws.Range("N77").Value = -73110
$ws.Range("H80").Value = 8002.6
$ws.Range("I80").Value = 7191.5835
$ws.Range("J80").Value = 9219.125
$ws.Range("K80").Value = 7191.5835
$ws.Range("L80").Value = 9219.125
$ws.Range("M80").Value = -6193.5835
$ws.Range("N80").Value = -11215.125
$ws.Range("H83").Value = 8002.6
$ws.Range("I83").Value = 7191.5835
$ws.Range("J83").Value = 9219.125
$ws.Range("K83").Value = 35957.9175
$ws.Range("L83").Value = 46095.625
$ws.Range("M83").Value = -30965.9175
$ws.Range("N83").Value = -56079.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 261.64285
$ws.Range("I55").Value = 233
$ws.Range("K55").Value = 233
$ws.Range("M55").Value = -60
$ws.Range("H93").Value = 929.3333
$ws.Range("I93").Value = 929.3333
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 929.3333
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 318.6667
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 6965.6665
$ws.Range("I122").Value = 5798.7
$ws.Range("K122").Value = 17396.1
$ws.Range("M122").Value = -14946.1
$ws.Range("H136").Value = 3080.8823
$ws.Range("I136").Value = 1920.2174
$ws.Range("J136").Value = 5507.727
$ws.Range("K136").Value = 5760.6522
$ws.Range("L136").Value = 16523.181
$ws.Range("M136").Value = -3210.6522
$ws.Range("N136").Value = -21623.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1732.3
$ws.Range("I81").Value = 1963.2858
$ws.Range("K81").Value = 3926.5716
$ws.Range("M81").Value = -2865.5716
$ws.Range("H84").Value = 1732.3
$ws.Range("I84").Value = 1963.2858
$ws.Range("K84").Value = 19632.858
$ws.Range("M84").Value = -14328.858
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H136").Value = 10484961
$ws.Range("I136").Value = 16837522
$ws.Range("K136").Value = 50512566
$ws.Range("M136").Value = -50510016

Write-Output "Applied updates to Tonberry_Profits sheets."
